# Add a new "# Aged Care Active Outbreaks (Weekly)" metric block to the
# Metrics table, by inserting 3 rows just above the existing
# "# Aged Care Active Outbreaks" rows (before old row 65) and filling
# them in with the new weekly metrics, then growing the table/autofilter
# to cover the 3 extra rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert three new blank rows above row 65 (pushes the old rows 65-71
# down to 68-74), inheriting the existing row's formatting (style s="2"
# on column C carries into the new rows, matching the target rows).
$ws.Range("A65:F65").EntireRow.Insert()
$ws.Range("A65:F65").EntireRow.Insert()
$ws.Range("A65:F65").EntireRow.Insert()

# Grow the table (and its autofilter) to include the 3 new rows plus the
# rows that shifted down, i.e. from A1:F71 to A1:F74.
$lo.Resize($ws.Range("A1:F74"))

# New row 65: # Aged Care Active Outbreaks (Weekly)
$ws.Range("A65").Value = "Aged Care"
$ws.Range("B65").Value = 60
$ws.Range("C65").Value = "# Aged Care Active Outbreaks (Weekly)"
$ws.Range("D65").Value = 640
$ws.Range("F65").Value = "X"

# New row 66: # Aged Care Active Outbreaks (Weekly) per 1M
$ws.Range("A66").Value = "Aged Care"
$ws.Range("B66").Value = 60
$ws.Range("C66").Value = "# Aged Care Active Outbreaks (Weekly) per 1M"
$ws.Range("D66").Value = 650
$ws.Range("F66").Value = "X"

# New row 67: % Aged Care Active Outbreaks (Weekly) Change
$ws.Range("A67").Value = "Aged Care"
$ws.Range("B67").Value = 60
$ws.Range("C67").Value = "% Aged Care Active Outbreaks (Weekly) Change"
$ws.Range("D67").Value = 660
$ws.Range("F67").Value = "X"

# Renumber the "D" sort column for all the rows that shifted down so the
# sequence (10 per row) stays contiguous through to the end of the table.
$ws.Range("D68").Value = 670
$ws.Range("D69").Value = 680
$ws.Range("D70").Value = 690
$ws.Range("D71").Value = 700
$ws.Range("D72").Value = 710
$ws.Range("D73").Value = 720
$ws.Range("D74").Value = 730

# Match the author's final on-screen selection/scroll state.
$ws.Range("F64:F67").Select()
